$d = $word.ActiveDocument

$replacements = @(
    @("2025-01-04 Saturday", "2025-01-05 Sunday"),
    @("22+67=", "96-10="),
    @("10+75=", "2+61="),
    @("52-37=", "7+57="),
    @("72-65=", "15+43="),
    @("8+2=", "3+71="),
    @("41-29=", "51-8="),
    @("12+3=", "48-29="),
    @("13+27=", "6+8="),
    @("86-58=", "57+15="),
    @("62-16=", "58-33="),
    @("37-6=", "77+21="),
    @("80-57=", "60-17="),
    @("88-70=", "21+35="),
    @("72-17=", "31+47="),
    @("5+28=", "29-16="),
    @("28+56=", "91-61="),
    @("80-73=", "29+64="),
    @("73-26=", "77-33="),
    @("69+3=", "12+57="),
    @("85-12=", "81-3="),
    @("85-14=", "46+52="),
    @("28+28=", "13+75="),
    @("99-77=", "28-2="),
    @("15+17=", "3+16="),
    @("71+2=", "56+35="),
    @("75-29=", "91-0="),
    @("9+31=", "34+16="),
    @("28+36=", "36-20="),
    @("92-68=", "2+36="),
    @("85-6=", "52-29="),
    @("25+70=", "86-60="),
    @("0+86=", "44-4="),
    @("20+79=", "84+0="),
    @("22-13=", "55+35="),
    @("0+71=", "90-60="),
    @("79-19=", "57-31="),
    @("5+29=", "31-25="),
    @("24-11=", "74-53="),
    @("45-11=", "88-80="),
    @("20+17=", "37+28="),
    @("6+53=", "50+3="),
    @("91-82=", "0+4="),
    @("73-7=", "73-70="),
    @("9+71=", "6+24="),
    @("24+75=", "60-56="),
    @("92-60=", "26+24="),
    @("75-65=", "87-23="),
    @("50+25=", "49-37="),
    @("98-80=", "18+27="),
    @("9+64=", "97-57="),
    @("4+63=", "39+18="),
    @("27+38=", "60+36="),
    @("24+1=", "40+32="),
    @("54-16=", "7+75="),
    @("86-27=", "52+28="),
    @("28-5=", "98-51="),
    @("53+21=", "74+13="),
    @("60-11=", "31+13="),
    @("11-5=", "83-11="),
    @("58-15=", "80-66="),
    @("31-27=", "60+21="),
    @("17+42=", "9+50="),
    @("72+25=", "26+51="),
    @("26-19=", "44-2="),
    @("12+37=", "59+2="),
    @("78-60=", "98-47="),
    @("2+76=", "84-70="),
    @("62-21=", "78-77="),
    @("73-55=", "63-1="),
    @("27-3=", "81-55="),
    @("42-36=", "82-44="),
    @("59-4=", "89-17="),
    @("94-37=", "73-38="),
    @("97-50=", "83-42="),
    @("81-25=", "38-9="),
    @("83-65=", "98-35="),
    @("18+32=", "96-14="),
    @("44-6=", "95-7="),
    @("21+13=", "69+19="),
    @("29+11=", "77-23="),
    @("11+71=", "68+8="),
    @("81+1=", "24+59="),
    @("90-67=", "10+72="),
    @("39-17=", "95-73="),
    @("71+12=", "7+27="),
    @("86-21=", "45+0="),
    @("66-10=", "99-45="),
    @("40+49=", "18+38="),
    @("38+5=", "22+75="),
    @("25+1=", "10+32="),
    @("0+82=", "78+13="),
    @("45-42=", "18+61="),
    @("50-42=", "17-10="),
    @("49+31=", "94-72="),
    @("87-51=", "97-27="),
    @("72-37=", "98-67="),
    @("78-19=", "95-73="),
    @("18+40=", "90-62="),
    @("79-7=", "47+5="),
    @("62-39=", "69+14=")
)

foreach ($pair in $replacements) {
    $null = $d.Content.Find.Execute($pair[0], $true, $true, $false, $false, $false, $true, 1, $false, $pair[1], 2)
}
